# Fruta / hortaliza, semanal
# Insert a new weekly record at row 806 on Sheet1, shifting the existing
# rows 806-834 down to 807-835 (dimension grows from A1:T834 to A1:T835).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 806; this shifts rows 806..834
# down to 807..835 and carries their values/formatting with them, exactly
# like pressing "Insert Sheet Rows" in Excel.
$ws.Rows.Item(806).Insert()

# Populate the newly-inserted row 806 with the new record's data.
$ws.Cells.Item(806, 1).Value2 = 6
$ws.Cells.Item(806, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(806, 3).Value2 = "Metropolitana"
$ws.Cells.Item(806, 4).Value2 = 45075
$ws.Cells.Item(806, 5).Value2 = 13
$ws.Cells.Item(806, 6).Value2 = "Fruta"
$ws.Cells.Item(806, 7).Value2 = 100103
$ws.Cells.Item(806, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(806, 9).Value2 = 100103002
$ws.Cells.Item(806, 10).Value2 = "Ciruela"
$ws.Cells.Item(806, 11).Value2 = "Angeleno"
$ws.Cells.Item(806, 12).Value2 = "Segunda"
$ws.Cells.Item(806, 13).Value2 = 15
$ws.Cells.Item(806, 14).Value2 = 160000
$ws.Cells.Item(806, 15).Value2 = 160000
$ws.Cells.Item(806, 16).Value2 = 160000
$ws.Cells.Item(806, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(806, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(806, 19).Value2 = 356
$ws.Cells.Item(806, 20).Value2 = 450
